$d = $word.ActiveDocument

# --- Step 1: remove the old "Meta description" paragraph that sits right
# after the H1 title at the top of the document. ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Step 2: insert a new bold paragraph ("Play Book of Spells Slot for
# Free - Review 2021") right before the closing "Prompt: ..." paragraph
# at the end of the document, matching the original doc's run layout
# (a leading empty run followed by the bold text run). ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $lastPara.Range.Start
$insertionRange = $d.Range($insertPos, $insertPos)
$newTitle = "Play Book of Spells Slot for Free - Review 2021"
$insertionRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>" + $newTitle + "</w:t></w:r></w:p>")

# The InsertXML call merges the new paragraph content into the following
# ("Prompt: ...") paragraph, so split it back into its own paragraph right
# after the inserted title text.
$splitPos = $insertPos + $newTitle.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphBefore()

# --- Step 3: replace the "Prompt: ..." paragraph's text with the old
# meta-description copy text, keeping its existing (italic) run
# formatting untouched. ---
$oldPromptText = "Prompt: Create a cartoon image featuring a happy Maya warrior with glasses to fit the theme of the game " + [char]34 + "Book of Spells" + [char]34 + ". The Maya warrior should be standing in front of a cauldron with spell books and magic objects surrounding him. The background should be an enchanted forest with stars and sparkles. The image should be colorful and convey the theme of the game. Use bright colors for the Maya warrior" + [char]39 + "s clothing and accessories, and make the cauldron and spell books stand out. The image should be eye-catching and visually striking to attract potential players to the game."
$newPromptText = "Explore the world of Book of Spells slot for free with our review! Discover its pros, cons, and similar games available to play online."

$d.Content.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 1, $false, $newPromptText, 2)
